$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks up front -- they will be re-added after the
# column insert below shifts their target cells two columns to the right.
$ws.Cells.Hyperlinks.Delete()

# Insert two new columns before the old column B ("survey data source" etc.),
# pushing every existing column from B onward two slots to the right.
$ws.Range("B1:C1").EntireColumn.Insert()

# Give the two new columns a width close to their siblings.
$ws.Range("B:C").ColumnWidth = 28.71

# --- New column headers (row 1) ---
$ws.Range("B1").Value = "different from others because #1"
$ws.Range("C1").Value = "different from others because #2"

# --- Brazil row (row 3) ---
$ws.Range("B3").Value = "computes small area population weights from census microdata rather than from a prepared file"
$ws.Range("C3").Value = "???"

# --- European Union row (row 5) ---
$ws.Range("B5").Value = "chooses colors based on statistically significant differences"
$ws.Range("C5").Value = "???"

# --- India row (row 7) ---
# Column D7 held the old rich-text cell ("demographic and health survey" +
# bold red "with prevR" on a second line) shifted over from the old B7.
# Split that into plain-text cells: the base sentence stays (now plain) in
# D7, "uses the prevR package" becomes the new B7, and a brand-new question
# becomes the new C7.
$ws.Range("D7").Value = "demographic and health survey"
$ws.Range("B7").Value = "uses the prevR package"
$ws.Range("C7").Value = "??? Do we also convert this to ggplot2?  Or is that too hard?"

# --- Re-create the hyperlinks at their shifted locations ---
$ws.Hyperlinks.Add($ws.Range("I5"), "http://epp.eurostat.ec.europa.eu/portal/page/portal/gisco_Geographical_information_maps/popups/references/administrative_units_statistical_units_1")
$ws.Hyperlinks.Add($ws.Range("G5"), "http://epp.eurostat.ec.europa.eu/portal/page/portal/population/data/database")
$ws.Hyperlinks.Add($ws.Range("I3"), "ftp://geoftp.ibge.gov.br/malhas_digitais/censo_2010/")

# Put the selection where the author's cursor ended up.
$ws.Range("C7").Select()
